$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 539, shifting existing rows 539:584 down
# to 541:586 (dimension grows from A1:R584 to A1:R586).
$ws.Rows("539:540").Insert()

# New row 539 (Huracan / Primera)
$ws.Range("A539").Value = 1
$ws.Range("B539").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C539").Value = "Arica y Parinacota"
$ws.Range("D539").Value = 45166
$ws.Range("E539").Value = 15
$ws.Range("F539").Value = 100112032
$ws.Range("G539").Value = "Zapallo italiano"
$ws.Range("H539").Value = "Huracán"
$ws.Range("I539").Value = "Primera"
$ws.Range("J539").Value = 140
$ws.Range("K539").Value = 7000
$ws.Range("L539").Value = 8000
$ws.Range("M539").Value = 7500
$ws.Range("N539").Value = "`$/caja 70 unidades"
$ws.Range("O539").Value = "Región de Arica y Parinacota"
$ws.Range("P539").Value = 107
$ws.Range("Q539").Value = 70
$ws.Range("R539").Value = "Hortaliza"

# New row 540 (Huracan / Segunda)
$ws.Range("A540").Value = 1
$ws.Range("B540").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C540").Value = "Arica y Parinacota"
$ws.Range("D540").Value = 45166
$ws.Range("E540").Value = 15
$ws.Range("F540").Value = 100112032
$ws.Range("G540").Value = "Zapallo italiano"
$ws.Range("H540").Value = "Huracán"
$ws.Range("I540").Value = "Segunda"
$ws.Range("J540").Value = 160
$ws.Range("K540").Value = 6000
$ws.Range("L540").Value = 7000
$ws.Range("M540").Value = 6500
$ws.Range("N540").Value = "`$/caja 100 unidades"
$ws.Range("O540").Value = "Región de Arica y Parinacota"
$ws.Range("P540").Value = 65
$ws.Range("Q540").Value = 100
$ws.Range("R540").Value = "Hortaliza"
